$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('253 Hoddle Street Apartment Complex Collingwood', 5, 'new'),
  @('3535 Opal Meadow Heights Aged Care Community Meadow Heights', 25, 'old'),
  @('3535 Opal Meadow Heights Aged Care Community Meadow Heights', 26, 'new'),
  @('7 Chefs Fawkner', 6, 'old'),
  @('Acquire BPO Southbank', 8, 'new'),
  @('Acquire BPO Southbank', 9, 'old'),
  @('Al Haj Halal Meats Glenroy', 43, 'new'),
  @('Al Haj Halal Meats Glenroy', 50, 'old'),
  @('CS Square Caroline Springs', 7, 'new'),
  @('CS Square Caroline Springs', 8, 'old'),
  @('Can Panel Cambellfield', 5, 'old'),
  @('Cannie Road Construction Site Cannie', 7, 'old'),
  @('Cannie Road Construction Site Cannie', 9, 'new'),
  @('Caroline Springs Police Station', 8, 'old'),
  @('Caroline Springs Police Station', 11, 'new'),
  @('City of Hobsons Bay Community', 8, 'new'),
  @('City of Hobsons Bay Community', 9, 'old'),
  @('City of Wyndham Community', 6, 'old'),
  @('City of Wyndham Community', 7, 'new'),
  @('Classy Cabinets and Kitchens Craigieburn', 11, 'new'),
  @('Classy Cabinets and Kitchens Craigieburn', 13, 'old'),
  @('Community Kids Meadow Heights', 16, 'new'),
  @('Community Kids Meadow Heights', 19, 'old'),
  @('Construction Site Olea Apartment Caulfield North', 8, 'old'),
  @('Construction Site Olea Apartment Caulfield North', 10, 'new'),
  @('Costco Wholesale Epping', 18, 'new'),
  @('Costco Wholesale Epping', 20, 'old'),
  @('Croydon Orthodontics', 5, 'old'),
  @('Direct Freight Express Cambellfield', 10, 'old'),
  @('Direct Freight Express Campbellfield', 12, 'new'),
  @('Fitzroy Community School Fitzroy North', 45, 'new'),
  @('Fitzroy Community School Fitzroy North', 49, 'old'),
  @('Fonterra Manufacturing Workplace Campbellfield', 5, 'new'),
  @('Glenroy West Primary School', 5, 'new'),
  @('Glenroy West Primary School', 6, 'old'),
  @('ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine', 7, 'new'),
  @('KFC Fawkner', 8, 'new'),
  @('KFC Fawkner', 9, 'old'),
  @('Kasr Sweets Coolaroo', 7, 'new'),
  @('Kasr Sweets Coolaroo', 9, 'old'),
  @('Learning Nest Early Learning Centre Meadow Heights', 8, 'new'),
  @('Learning Nest Early Learning Centre Meadow Heights', 9, 'old'),
  @('Lineage Logistics Laverton North', 7, 'old'),
  @('Lineage Logistics Laverton North', 8, 'new'),
  @('Linfox Somerton National Distribution Centre Somerton', 6, 'old'),
  @('Linfox Somerton National Distribution Centre Somerton', 9, 'new'),
  @('Melbourne Metropolitan Remand Centre Ravenhall', 9, 'old'),
  @('Melbourne Metropolitan Remand Centre Ravenhall', 10, 'new'),
  @('Melbourne Truck Repairs Campbellfield', 6, 'old'),
  @('Melbourne Truck Repairs Campbellfield', 8, 'new'),
  @('Melbourne Youth Justice Centre Parkville', 5, 'old'),
  @('Melbourne Youth Justice Centre Parkville', 6, 'new'),
  @('Mercy Hospital for Women Heidelberg', 5, 'new'),
  @('Mercy Hospital for Women Heidelberg', 6, 'old'),
  @('Newbury Child and Community Centre Craigieburn', 9, 'new'),
  @('Newbury Child and Community Centre Craigieburn', 10, 'old'),
  @('Newport Football Club Altona North', 5, 'old'),
  @('Nido Early School Moonee Ponds', 16, 'new'),
  @('Nido Early School Moonee Ponds', 18, 'old'),
  @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 46, 'new'),
  @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 52, 'old'),
  @('OnQ Plumbing and Excavations Craigieburn', 15, 'old'),
  @('OnQ Plumbing and Excavations Craigieburn', 18, 'new'),
  @('Oporto Coolaroo', 15, 'new'),
  @('Oporto Coolaroo', 16, 'old'),
  @('Oscar Romero Catholic Primary School Craigieburn', 5, 'new'),
  @('Panorama Construction Site Whitehorse Rd Box Hill', 20, 'new'),
  @('Panorama Construction Site Whitehorse Rd Box Hill', 27, 'old'),
  @('Serco Mill Park', 5, 'new'),
  @('Tek Foods Somerton', 11, 'old'),
  @('Tek Foods Somerton', 13, 'new'),
  @('The Homestead Child and Family Centre Roxburgh Park', 17, 'new'),
  @('The Homestead Child and Family Centre Roxburgh Park', 18, 'old'),
  @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 8, 'old'),
  @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 9, 'new'),
  @('ThorwestenCabinets Pakenham', 12, 'new'),
  @('Tip Top Warehouse Dandenong', 5, 'new'),
  @('Tip Top Warehouse Dandenong', 8, 'old'),
  @('Total Window Concepts Hoppers Crossing', 6, 'new'),
  @('Yara Childcare Centre Truganina', 10, 'new'),
  @('Yara Childcare Centre Truganina', 13, 'old')
)

$row = 2
foreach ($item in $data) {
  $ws.Cells.Item($row, 1).Value = $item[0]
  $ws.Cells.Item($row, 2).Value = $item[1]
  $ws.Cells.Item($row, 3).Value = $item[2]
  $row = $row + 1
}
